$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Connector CN-018 (row 76-78): pin 3 signal name change "+9-13.2V " -> "12V AUX"
# (leading apostrophe preserves the cell's existing quote-prefix style)
$ws.Range("C78").Value = "'12V AUX"

# Connector CN-019 (row 80-82): populate pin # and signal name columns
$ws.Range("B80").Value = 1
$ws.Range("C80").Value = "GND"

$ws.Range("B81").Value = 2
$ws.Range("C81").Value = "Signal "

$ws.Range("B82").Value = 3
$ws.Range("C82").Value = "12V AUX"

# Update the active selection to match the saved view state
$ws.Range("C83").Select() | Out-Null
